# Reorder worksheet tabs: move "Npc" so it becomes the first sheet
# (before "디스크립션"), then leave "디스크립션" as the active/selected sheet.
$wb = $excel.ActiveWorkbook

$npc = $wb.Worksheets.Item("Npc")
$npc.Move($wb.Worksheets.Item(1))

$wb.Worksheets.Item("디스크립션").Activate()
